$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: move value from D2 to C2 (new meanEMG value)
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 42.093751076109172

# Row 3: clear legmaxROM values in B3 and C3
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update the selection to match the new active range
[void]$ws.Range("B1:E3").Select()
